$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("metodo_directo_pbi")
$ws1.Cells.Item(2, 2).Value = 0.009990633367597335
$ws1.Cells.Item(2, 3).Value = 0.004871960842303819
$ws1.Cells.Item(3, 2).Value = 0.0169396038823721
$ws1.Cells.Item(3, 3).Value = 0.00799630893428064
$ws1.Cells.Item(4, 2).Value = 0.0286009585862577
$ws1.Cells.Item(4, 3).Value = 0.01389000088686124
$ws1.Cells.Item(5, 2).Value = 0.0158742471577314
$ws1.Cells.Item(5, 3).Value = 0.007742028046995103
$ws1.Cells.Item(6, 2).Value = -0.01523173853828302
$ws1.Cells.Item(6, 3).Value = -0.007727061243745569
$ws1.Cells.Item(7, 2).Value = 0.00005900023827991279
$ws1.Cells.Item(7, 3).Value = 0.00003135809694375898
$ws1.Cells.Item(8, 2).Value = 0.02091530605131297
$ws1.Cells.Item(8, 3).Value = 0.01081517322145661
$ws1.Cells.Item(9, 2).Value = 0.0513187319045086
$ws1.Cells.Item(9, 3).Value = 0.02635750221095647
$ws1.Cells.Item(10, 2).Value = 0.06574324644245021
$ws1.Cells.Item(10, 3).Value = 0.03343892775782109
$ws1.Cells.Item(11, 2).Value = 0.06989161493457444
$ws1.Cells.Item(11, 3).Value = 0.03533884675085421
$ws1.Cells.Item(12, 2).Value = 0.1600578399035218
$ws1.Cells.Item(12, 3).Value = 0.08100359237382432
$ws1.Cells.Item(13, 2).Value = 0.09111167441007032
$ws1.Cells.Item(13, 3).Value = 0.03984572516573932
$ws1.Cells.Item(14, 2).Value = 0.116062682203614
$ws1.Cells.Item(14, 3).Value = 0.05402523911473674
$ws1.Cells.Item(15, 2).Value = 0.1284670746552036
$ws1.Cells.Item(15, 3).Value = 0.05918944848216053
$ws1.Cells.Item(16, 2).Value = 0.1493599513095874
$ws1.Cells.Item(16, 3).Value = 0.06446795251740316
$ws1.Cells.Item(17, 2).Value = 0.1023046688651035
$ws1.Cells.Item(17, 3).Value = 0.04259304184329803
$ws1.Cells.Item(18, 2).Value = 0.09503219858379991
$ws1.Cells.Item(18, 3).Value = 0.04173256491423031
$ws1.Cells.Item(19, 2).Value = 0.01321584916741356
$ws1.Cells.Item(19, 3).Value = 0.005685224386666409
$ws1.Cells.Item(20, 2).Value = -0.02218241598321194
$ws1.Cells.Item(20, 3).Value = -0.009738806156855918
$ws1.Cells.Item(21, 2).Value = 0.02235913486482052
$ws1.Cells.Item(21, 3).Value = 0.009734289135004641
$ws1.Cells.Item(22, 2).Value = 0.04196129166117218
$ws1.Cells.Item(22, 3).Value = 0.01942539529351357

$ws2 = $wb.Worksheets.Item("metodo_indirecto_pbi")
$ws2.Cells.Item(54, 2).Value = -0.004901212763098724
$ws2.Cells.Item(54, 3).Value = -0.003246665888557902
$ws2.Cells.Item(55, 2).Value = 0.009387308803493815
$ws2.Cells.Item(55, 3).Value = 0.006279023819487015
$ws2.Cells.Item(56, 2).Value = 0.0129795170975493
$ws2.Cells.Item(56, 3).Value = 0.007992177513665126
$ws2.Cells.Item(57, 2).Value = 0.03761076148538569
$ws2.Cells.Item(57, 3).Value = 0.02060566025693597
$ws2.Cells.Item(58, 2).Value = -0.007920198641727725
$ws2.Cells.Item(58, 3).Value = -0.004419107826238254
$ws2.Cells.Item(59, 2).Value = -0.02062552237906596
$ws2.Cells.Item(59, 3).Value = -0.01161617805401889
$ws2.Cells.Item(60, 2).Value = -0.0063765978523947
$ws2.Cells.Item(60, 3).Value = -0.003912175993438178
$ws2.Cells.Item(61, 2).Value = 0.007097576608135333
$ws2.Cells.Item(61, 3).Value = 0.004025055359196407
$ws2.Cells.Item(62, 2).Value = 0.01734149166956163
$ws2.Cells.Item(62, 3).Value = 0.009934527433127527
$ws2.Cells.Item(63, 2).Value = 0.01699483601883259
$ws2.Cells.Item(63, 3).Value = 0.01074326449997213
$ws2.Cells.Item(64, 2).Value = 0.03765134207121663
$ws2.Cells.Item(64, 3).Value = 0.0239204711384449
$ws2.Cells.Item(65, 2).Value = 0.08493822413534784
$ws2.Cells.Item(65, 3).Value = 0.05085864499326751
$ws2.Cells.Item(66, 2).Value = 0.08836719139154742
$ws2.Cells.Item(66, 3).Value = 0.05661853707127002
$ws2.Cells.Item(67, 2).Value = 0.05465840093706187
$ws2.Cells.Item(67, 3).Value = 0.04793361363946315
$ws2.Cells.Item(68, 2).Value = 0.05719544102851242
$ws2.Cells.Item(68, 3).Value = 0.04824451248738074
$ws2.Cells.Item(69, 2).Value = 0.06398061523838239
$ws2.Cells.Item(69, 3).Value = 0.04703370451527033
$ws2.Cells.Item(70, 2).Value = 0.141214599277904
$ws2.Cells.Item(70, 3).Value = 0.09154275765902263
$ws2.Cells.Item(71, 2).Value = 0.1652441488932682
$ws2.Cells.Item(71, 3).Value = 0.0966575387170544
$ws2.Cells.Item(72, 2).Value = 0.1531734326065899
$ws2.Cells.Item(72, 3).Value = 0.1007053600286558
$ws2.Cells.Item(73, 2).Value = 0.1034819358691727
$ws2.Cells.Item(73, 3).Value = 0.08867874158916557
$ws2.Cells.Item(74, 2).Value = 0.07314611378719588
$ws2.Cells.Item(74, 3).Value = 0.06042928722933458
$ws2.Cells.Item(75, 2).Value = 0.1075605416680581
$ws2.Cells.Item(75, 3).Value = 0.06247441742099008
$ws2.Cells.Item(76, 2).Value = 0.1110972247952469
$ws2.Cells.Item(76, 3).Value = 0.05861640702145614
$ws2.Cells.Item(77, 2).Value = 0.07049759290623234
$ws2.Cells.Item(77, 3).Value = 0.0337956883315966
$ws2.Cells.Item(78, 2).Value = 0.05597201346912049
$ws2.Cells.Item(78, 3).Value = 0.02930917591928976
$ws2.Cells.Item(79, 2).Value = 0.0171520373095995
$ws2.Cells.Item(79, 3).Value = 0.01140100880923048
$ws2.Cells.Item(80, 2).Value = 0.07880279842802054
$ws2.Cells.Item(80, 3).Value = 0.05340293123463417
$ws2.Cells.Item(81, 2).Value = 0.1092788477922576
$ws2.Cells.Item(81, 3).Value = 0.05498605490750465
$ws2.Cells.Item(82, 2).Value = 0.1096926821680189
$ws2.Cells.Item(82, 3).Value = 0.05536083519140974
$ws2.Cells.Item(83, 2).Value = 0.07789626916676194
$ws2.Cells.Item(83, 3).Value = 0.03459952813125559
$ws2.Cells.Item(84, 2).Value = 0.0520423435993547
$ws2.Cells.Item(84, 3).Value = 0.01839028920722597
$ws2.Cells.Item(85, 2).Value = 0.03641018670097225
$ws2.Cells.Item(85, 3).Value = 0.01512391535739226
$ws2.Cells.Item(86, 2).Value = 0.04653681169069657
$ws2.Cells.Item(86, 3).Value = 0.01938601454943325
$ws2.Cells.Item(87, 2).Value = 0.04372915311954686
$ws2.Cells.Item(87, 3).Value = 0.02035523739401826
$ws2.Cells.Item(88, 2).Value = 0.03286385994041376
$ws2.Cells.Item(88, 3).Value = 0.01614851990478608
$ws2.Cells.Item(89, 2).Value = 0.01883044223656507
$ws2.Cells.Item(89, 3).Value = 0.009182718837161367
$ws2.Cells.Item(90, 2).Value = 0.02704993131650211
$ws2.Cells.Item(90, 3).Value = 0.01276887045056065
$ws2.Cells.Item(91, 2).Value = 0.05375678758173594
$ws2.Cells.Item(91, 3).Value = 0.02610688117089517
$ws2.Cells.Item(92, 2).Value = 0.04282621701312624
$ws2.Cells.Item(92, 3).Value = 0.0208867715090878
$ws2.Cells.Item(93, 2).Value = 0.04769463717964934
$ws2.Cells.Item(93, 3).Value = 0.02419549032824443
$ws2.Cells.Item(94, 2).Value = 0.03520371409022449
$ws2.Cells.Item(94, 3).Value = 0.0187104579812769
$ws2.Cells.Item(95, 2).Value = 0.04438140761208084
$ws2.Cells.Item(95, 3).Value = 0.0229493467587387
$ws2.Cells.Item(96, 2).Value = 0.06919930681393717
$ws2.Cells.Item(96, 3).Value = 0.03554103569314353
$ws2.Cells.Item(97, 2).Value = 0.08798484129787348
$ws2.Cells.Item(97, 3).Value = 0.04475164995872836
$ws2.Cells.Item(98, 2).Value = 0.09037406123221664
$ws2.Cells.Item(98, 3).Value = 0.04569525404624369
$ws2.Cells.Item(99, 2).Value = 0.1926888178556288
$ws2.Cells.Item(99, 3).Value = 0.09751778773210845
$ws2.Cells.Item(100, 2).Value = 0.1213655162815073
$ws2.Cells.Item(100, 3).Value = 0.05307659021373995
$ws2.Cells.Item(101, 2).Value = 0.1396652742196356
$ws2.Cells.Item(101, 3).Value = 0.06501185129001037
$ws2.Cells.Item(102, 2).Value = 0.1569866918309726
$ws2.Cells.Item(102, 3).Value = 0.07232947222821964
$ws2.Cells.Item(103, 2).Value = 0.1815875866490965
$ws2.Cells.Item(103, 3).Value = 0.07837830563816162
$ws2.Cells.Item(104, 2).Value = 0.1465107452040416
$ws2.Cells.Item(104, 3).Value = 0.06099759053222571
$ws2.Cells.Item(105, 2).Value = 0.1288625255172258
$ws2.Cells.Item(105, 3).Value = 0.05658885926349632
$ws2.Cells.Item(106, 2).Value = 0.05002668700437171
$ws2.Cells.Item(106, 3).Value = 0.02152059525941483
$ws2.Cells.Item(107, 2).Value = 0.0198697305366796
$ws2.Cells.Item(107, 3).Value = 0.008723461602745768
$ws2.Cells.Item(108, 2).Value = 0.03037512153385279
$ws2.Cells.Item(108, 3).Value = 0.01322413489202782
$ws2.Cells.Item(109, 2).Value = 0.02864378697744449
$ws2.Cells.Item(109, 3).Value = 0.01326024206387626
